$d = $word.ActiveDocument

function New-XmlPkg([string]$bodyXml) {
    return '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml" pkg:padding="512">' +
           '<pkg:xmlData>' +
           '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
           '<w:body>' + $bodyXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
}

# ---------------------------------------------------------------------------
# Locate the "v1.1.3 change log (...) (TBA)" paragraph.
# ---------------------------------------------------------------------------
$count = $d.Paragraphs.Count
$p127Index = -1
for ($i = 1; $i -le $count; $i++) {
    if ($d.Paragraphs($i).Range.Text -like "*start=XX8*TBA*") {
        $p127Index = $i
        break
    }
}
Write-Output "p127Index=$p127Index"

# ---------------------------------------------------------------------------
# Step 1: merge the trailing runs of that paragraph:
#   "ral/source/list?start=XX8&num=XX" + ") (" + "TBA" + ")"
#   -> single run "ral/source/list?start=XX8&num=XX) (TBA)"
# ---------------------------------------------------------------------------
$p127 = $d.Paragraphs($p127Index)
$bodyXml = '<w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>' +
           '<w:proofErr w:type="gramStart"/>' +
           '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>v1.1.3</w:t></w:r>' +
           '<w:proofErr w:type="gramEnd"/>' +
           '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> change log (</w:t></w:r>' +
           '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>http://code.google.com/p/subcent</w:t></w:r>' +
           '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>ral/source/list?start=XX8&amp;num=XX) (TBA)</w:t></w:r>' +
           '</w:p>'
$p127.Range.InsertXML((New-XmlPkg $bodyXml))
Write-Output "afterStep1=[$($d.Paragraphs($p127Index).Range.Text)]"

# ---------------------------------------------------------------------------
# Step 2: insert two new "ListParagraph" bullet paragraphs right before the
# "Added Icelandic translation (...)" paragraph, for the SubsCenter.org and
# Sratim.co.il changelog entries.
# ---------------------------------------------------------------------------
$icelandicIndex = $p127Index + 1
$icelandicPara = $d.Paragraphs($icelandicIndex)
Write-Output "icelandicBefore=[$($icelandicPara.Range.Text)]"

# Insert two blank paragraphs before it; each inherits the ListParagraph /
# numId=14 paragraph formatting from the paragraph it is inserted before.
$icelandicPara.Range.InsertParagraphBefore()
$icelandicPara = $d.Paragraphs($icelandicIndex + 1)
$icelandicPara.Range.InsertParagraphBefore()

$subsCenterIndex = $icelandicIndex
$sratimIndex = $icelandicIndex + 1
$icelandicIndex = $icelandicIndex + 2

Write-Output "subsCenterIndex=$subsCenterIndex sratimIndex=$sratimIndex icelandicIndex=$icelandicIndex"

$pPrListParagraph = '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="14"/></w:numPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>'

# --- "Updated SubsCenter.org (yoavain)" -------------------------------------
$bodyXml = '<w:p>' + $pPrListParagraph +
           '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Updated </w:t></w:r>' +
           '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>SubsCenter.org</w:t></w:r>' +
           '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> (</w:t></w:r>' +
           '<w:proofErr w:type="spellStart"/>' +
           '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>yoavain</w:t></w:r>' +
           '<w:proofErr w:type="spellEnd"/>' +
           '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>)</w:t></w:r>' +
           '</w:p>'
$d.Paragraphs($subsCenterIndex).Range.InsertXML((New-XmlPkg $bodyXml))
Write-Output "afterSubsCenter=[$($d.Paragraphs($subsCenterIndex).Range.Text)]"

# --- "Added Sratim.co.il subtitle provider (SubtitleDownloader) (yoavain)" --
$bodyXml = '<w:p>' + $pPrListParagraph +
           '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Added Sratim.co.il</w:t></w:r>' +
           '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' +
           '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>subtitl</w:t></w:r>' +
           '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>e provider (</w:t></w:r>' +
           '<w:proofErr w:type="spellStart"/>' +
           '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>SubtitleDownloader</w:t></w:r>' +
           '<w:proofErr w:type="spellEnd"/>' +
           '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>)</w:t></w:r>' +
           '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> (</w:t></w:r>' +
           '<w:proofErr w:type="spellStart"/>' +
           '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>yoavain</w:t></w:r>' +
           '<w:proofErr w:type="spellEnd"/>' +
           '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>)</w:t></w:r>' +
           '</w:p>'
$d.Paragraphs($sratimIndex).Range.InsertXML((New-XmlPkg $bodyXml))
Write-Output "afterSratim=[$($d.Paragraphs($sratimIndex).Range.Text)]"

# --- merge "Added " + "Icelandic " + "translation (" into one run ----------
$bodyXml = '<w:p>' + $pPrListParagraph +
           '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Added Icelandic translation (</w:t></w:r>' +
           '<w:proofErr w:type="spellStart"/>' +
           '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>kjarrig</w:t></w:r>' +
           '<w:proofErr w:type="spellEnd"/>' +
           '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>)</w:t></w:r>' +
           '</w:p>'
$d.Paragraphs($icelandicIndex).Range.InsertXML((New-XmlPkg $bodyXml))
Write-Output "afterIcelandic=[$($d.Paragraphs($icelandicIndex).Range.Text)]"

for ($i = $p127Index; $i -le $icelandicIndex + 2; $i++) {
    Write-Output "chk $i [$($d.Paragraphs($i).Range.Text)]"
}
